# "If sheet name but not xls path specified, assume sheet is in same file"
#
# On the "Chain List" sheet, the ChainFile column (D) previously repeated the
# workbook's own path ("excelData/cementFactory.xlsx") for every chain. Now
# that a missing file path is understood to mean "this workbook", the sample
# data is updated to show that: the cement chain keeps a descriptive
# placeholder ("here"), while the other two chains use the new shorthand
# values ("thisfile" / "same") instead of repeating the literal path.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chain List")

# Update the ChainFile values (column D) for rows 2-4. Order matters so the
# shared-string table ends up reused/appended the way Excel would do it.
$ws.Range("D2").Value = "here"
$ws.Range("D3").Value = "thisfile"
$ws.Range("D4").Value = "same"

# Column D no longer needs to be wide enough to hold a full relative path;
# shrink it back down to match column B's width.
$ws.Columns.Item(4).ColumnWidth = 11

# Move the active selection (the workbook was left scrolled/selected at F10).
$ws.Range("F10").Select() | Out-Null

# Window was shifted slightly on screen.
$win = $excel.ActiveWindow
$win.Left = 1160
$win.Top = 460
